$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for rule R10 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Move the active selection to E8, matching the saved view state
$ws.Range("E8").Select()
